$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 222; this shifts all existing rows 222..334
# down to 223..335 (and the sheet dimension grows from A1:T334 to A1:T335).
$ws.Range("A222").EntireRow.Insert()

# Populate the newly inserted row 222 with the new pricing record.
$ws.Range("A222").Value = 7
$ws.Range("B222").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C222").Value = 'Ñuble'
$ws.Range("D222").Value = 45089
$ws.Range("E222").Value = 16
$ws.Range("F222").Value = 'Fruta'
$ws.Range("G222").Value = 100108
$ws.Range("H222").Value = 'Tropicales y subtropicales'
$ws.Range("I222").Value = 100108005
$ws.Range("J222").Value = 'Piña'
$ws.Range("K222").Value = 'Caramelo'
$ws.Range("L222").Value = 'Primera'
$ws.Range("M222").Value = 40
$ws.Range("N222").Value = 23000
$ws.Range("O222").Value = 23000
$ws.Range("P222").Value = 23000
$ws.Range("Q222").Value = '$/caja 12 unidades'
$ws.Range("R222").Value = 'Ecuador'
$ws.Range("S222").Value = 1917
$ws.Range("T222").Value = 12
